# Auto-generated edit script: updates Leve profit-calc cells (H:N) across 8 crafting sheets
# per the scheduled-runner price refresh (currentAveragePrice* / LevePrice* / LeveProfit*).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(70, 8).Value = 243327
$ws.Cells.Item(70, 9).Value = 283500
$ws.Cells.Item(70, 11).Value = 850500
$ws.Cells.Item(70, 13).Value = -850230

$ws.Cells.Item(73, 8).Value = 243327
$ws.Cells.Item(73, 9).Value = 283500
$ws.Cells.Item(73, 11).Value = 850500
$ws.Cells.Item(73, 13).Value = -849564

$ws.Cells.Item(80, 8).Value = 1051.1
$ws.Cells.Item(80, 10).Value = 1494.25
$ws.Cells.Item(80, 12).Value = 4482.75
$ws.Cells.Item(80, 14).Value = -6478.75

$ws.Cells.Item(83, 8).Value = 1051.1
$ws.Cells.Item(83, 10).Value = 1494.25
$ws.Cells.Item(83, 12).Value = 13448.25
$ws.Cells.Item(83, 14).Value = -23432.25

$ws.Cells.Item(135, 8).Value = 949.5238000000001
$ws.Cells.Item(135, 9).Value = 1033.7894
$ws.Cells.Item(135, 10).Value = 149
$ws.Cells.Item(135, 11).Value = 9304.104599999999
$ws.Cells.Item(135, 12).Value = 1341
$ws.Cells.Item(135, 13).Value = -6769.104599999999
$ws.Cells.Item(135, 14).Value = -6411

$ws.Cells.Item(137, 8).Value = 3606.5483
$ws.Cells.Item(137, 10).Value = 5046.643
$ws.Cells.Item(137, 12).Value = 15139.929
$ws.Cells.Item(137, 14).Value = -20239.929

$ws.Cells.Item(138, 8).Value = 2818.2192
$ws.Cells.Item(138, 10).Value = 2944.5557
$ws.Cells.Item(138, 12).Value = 8833.667099999999
$ws.Cells.Item(138, 14).Value = -19113.6671

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 3635.5557
$ws.Cells.Item(61, 10).Value = 4186.381
$ws.Cells.Item(61, 12).Value = 4186.381
$ws.Cells.Item(61, 14).Value = -4610.381

$ws.Cells.Item(63, 8).Value = 2658.1
$ws.Cells.Item(63, 9).Value = 2602.465
$ws.Cells.Item(63, 11).Value = 2602.465
$ws.Cells.Item(63, 13).Value = -1916.465

$ws.Cells.Item(66, 8).Value = 2658.1
$ws.Cells.Item(66, 9).Value = 2602.465
$ws.Cells.Item(66, 11).Value = 13012.325
$ws.Cells.Item(66, 13).Value = -9580.325000000001

$ws.Cells.Item(102, 8).Value = 3749.6924
$ws.Cells.Item(102, 9).Value = 3749.6924
$ws.Cells.Item(102, 11).Value = 3749.6924
$ws.Cells.Item(102, 13).Value = -2127.6924

$ws.Cells.Item(132, 8).Value = 2952.6843
$ws.Cells.Item(132, 9).Value = 3109.5454
$ws.Cells.Item(132, 10).Value = 2737
$ws.Cells.Item(132, 11).Value = 9328.636200000001
$ws.Cells.Item(132, 12).Value = 8211
$ws.Cells.Item(132, 13).Value = -6798.636200000001
$ws.Cells.Item(132, 14).Value = -13271

$ws.Cells.Item(135, 8).Value = 67500
$ws.Cells.Item(135, 10).Value = 67500
$ws.Cells.Item(135, 12).Value = 67500
$ws.Cells.Item(135, 14).Value = -77640

$ws.Cells.Item(136, 8).Value = 3635.5557
$ws.Cells.Item(136, 10).Value = 4186.381
$ws.Cells.Item(136, 12).Value = 12559.143
$ws.Cells.Item(136, 14).Value = -17659.143

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(53, 8).Value = 55700
$ws.Cells.Item(53, 9).Value = 55700
$ws.Cells.Item(53, 11).Value = 55700
$ws.Cells.Item(53, 13).Value = -55126

$ws.Cells.Item(80, 8).Value = 1896.2307
$ws.Cells.Item(80, 9).Value = 370
$ws.Cells.Item(80, 10).Value = 2354.1
$ws.Cells.Item(80, 11).Value = 370
$ws.Cells.Item(80, 12).Value = 2354.1
$ws.Cells.Item(80, 13).Value = 628
$ws.Cells.Item(80, 14).Value = -4350.1

$ws.Cells.Item(82, 8).Value = 14500.728
$ws.Cells.Item(82, 9).Value = 11950.8
$ws.Cells.Item(82, 11).Value = 11950.8
$ws.Cells.Item(82, 13).Value = -11567.8

$ws.Cells.Item(83, 8).Value = 1896.2307
$ws.Cells.Item(83, 9).Value = 370
$ws.Cells.Item(83, 10).Value = 2354.1
$ws.Cells.Item(83, 11).Value = 1850
$ws.Cells.Item(83, 12).Value = 11770.5
$ws.Cells.Item(83, 13).Value = 3142
$ws.Cells.Item(83, 14).Value = -21754.5

$ws.Cells.Item(85, 8).Value = 14500.728
$ws.Cells.Item(85, 9).Value = 11950.8
$ws.Cells.Item(85, 11).Value = 11950.8
$ws.Cells.Item(85, 13).Value = -10624.8

$ws.Cells.Item(94, 8).Value = 4792.75
$ws.Cells.Item(94, 9).Value = 3468.6
$ws.Cells.Item(94, 10).Value = 6999.6665
$ws.Cells.Item(94, 11).Value = 3468.6
$ws.Cells.Item(94, 12).Value = 6999.6665
$ws.Cells.Item(94, 13).Value = -3017.6
$ws.Cells.Item(94, 14).Value = -7901.6665

$ws.Cells.Item(105, 8).Value = 4764.1816
$ws.Cells.Item(105, 9).Value = 5978
$ws.Cells.Item(105, 11).Value = 5978
$ws.Cells.Item(105, 13).Value = -4231

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 3592.05
$ws.Cells.Item(58, 9).Value = 3518.9443
$ws.Cells.Item(58, 10).Value = 4250
$ws.Cells.Item(58, 11).Value = 3518.9443
$ws.Cells.Item(58, 12).Value = 4250
$ws.Cells.Item(58, 13).Value = -3315.9443
$ws.Cells.Item(58, 14).Value = -4656

$ws.Cells.Item(62, 8).Value = 11187.333
$ws.Cells.Item(62, 10).Value = 5195
$ws.Cells.Item(62, 12).Value = 5195
$ws.Cells.Item(62, 14).Value = -6443

$ws.Cells.Item(65, 8).Value = 11187.333
$ws.Cells.Item(65, 10).Value = 5195
$ws.Cells.Item(65, 12).Value = 25975
$ws.Cells.Item(65, 14).Value = -32215

$ws.Cells.Item(132, 8).Value = 3262.6365
$ws.Cells.Item(132, 9).Value = 3072.25
$ws.Cells.Item(132, 10).Value = 3770.3333
$ws.Cells.Item(132, 11).Value = 9216.75
$ws.Cells.Item(132, 12).Value = 11310.9999
$ws.Cells.Item(132, 13).Value = -6686.75
$ws.Cells.Item(132, 14).Value = -16370.9999

$ws.Cells.Item(133, 8).Value = 100326
$ws.Cells.Item(133, 10).Value = 100326
$ws.Cells.Item(133, 12).Value = 100326
$ws.Cells.Item(133, 14).Value = -105386

$ws.Cells.Item(136, 8).Value = 3592.05
$ws.Cells.Item(136, 9).Value = 3518.9443
$ws.Cells.Item(136, 10).Value = 4250
$ws.Cells.Item(136, 11).Value = 10556.8329
$ws.Cells.Item(136, 12).Value = 12750
$ws.Cells.Item(136, 13).Value = -8006.832900000001
$ws.Cells.Item(136, 14).Value = -17850

$ws.Cells.Item(141, 8).Value = 127773.39
$ws.Cells.Item(141, 10).Value = 127773.39
$ws.Cells.Item(141, 12).Value = 127773.39
$ws.Cells.Item(141, 14).Value = -138133.39

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(107, 8).Value = 716.8889
$ws.Cells.Item(107, 9).Value = 363.14285
$ws.Cells.Item(107, 10).Value = 942
$ws.Cells.Item(107, 11).Value = 1089.42855
$ws.Cells.Item(107, 12).Value = 2826
$ws.Cells.Item(107, 13).Value = 830.5714499999999
$ws.Cells.Item(107, 14).Value = -6666

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 6071.3335
$ws.Cells.Item(70, 9).Value = 5487.222
$ws.Cells.Item(70, 11).Value = 5487.222
$ws.Cells.Item(70, 13).Value = -5217.222

$ws.Cells.Item(73, 8).Value = 6071.3335
$ws.Cells.Item(73, 9).Value = 5487.222
$ws.Cells.Item(73, 11).Value = 5487.222
$ws.Cells.Item(73, 13).Value = -4551.222

$ws.Cells.Item(106, 8).Value = 29998.5
$ws.Cells.Item(106, 10).Value = 29998.5
$ws.Cells.Item(106, 12).Value = 29998.5
$ws.Cells.Item(106, 14).Value = -32522.5

$ws.Cells.Item(122, 8).Value = 2373.2856
$ws.Cells.Item(122, 9).Value = 1564.7273
$ws.Cells.Item(122, 10).Value = 3262.7
$ws.Cells.Item(122, 11).Value = 4694.1819
$ws.Cells.Item(122, 12).Value = 9788.099999999999
$ws.Cells.Item(122, 13).Value = -2244.1819
$ws.Cells.Item(122, 14).Value = -14688.1

$ws.Cells.Item(126, 8).Value = 3223.8113
$ws.Cells.Item(126, 9).Value = 2815.2307
$ws.Cells.Item(126, 10).Value = 3617.2593
$ws.Cells.Item(126, 11).Value = 8445.6921
$ws.Cells.Item(126, 12).Value = 10851.7779
$ws.Cells.Item(126, 13).Value = -5975.6921
$ws.Cells.Item(126, 14).Value = -15791.7779

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(25, 8).Value = 8718.25
$ws.Cells.Item(25, 9).Value = 7009
$ws.Cells.Item(25, 10).Value = 9288
$ws.Cells.Item(25, 11).Value = 7009
$ws.Cells.Item(25, 12).Value = 9288
$ws.Cells.Item(25, 13).Value = -6779
$ws.Cells.Item(25, 14).Value = -9748

$ws.Cells.Item(46, 8).Value = 1799.2059
$ws.Cells.Item(46, 9).Value = 1262.3334
$ws.Cells.Item(46, 11).Value = 1262.3334
$ws.Cells.Item(46, 13).Value = -1074.3334

$ws.Cells.Item(68, 8).Value = 2399.6
$ws.Cells.Item(68, 9).Value = 1666.3334
$ws.Cells.Item(68, 11).Value = 1666.3334
$ws.Cells.Item(68, 13).Value = -917.3334

$ws.Cells.Item(71, 8).Value = 2399.6
$ws.Cells.Item(71, 9).Value = 1666.3334
$ws.Cells.Item(71, 11).Value = 8331.666999999999
$ws.Cells.Item(71, 13).Value = -4587.666999999999

$ws.Cells.Item(100, 8).Value = 3924.5789
$ws.Cells.Item(100, 9).Value = 4237.4
$ws.Cells.Item(100, 11).Value = 4237.4
$ws.Cells.Item(100, 13).Value = -3696.4

$ws.Cells.Item(136, 8).Value = 4521.143
$ws.Cells.Item(136, 10).Value = 5000
$ws.Cells.Item(136, 12).Value = 15000
$ws.Cells.Item(136, 14).Value = -20100

$ws.Cells.Item(140, 8).Value = 66250
$ws.Cells.Item(140, 10).Value = 66250
$ws.Cells.Item(140, 12).Value = 66250
$ws.Cells.Item(140, 14).Value = -76610

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(41, 8).Value = 9154.25
$ws.Cells.Item(41, 9).Value = 11979.4
$ws.Cells.Item(41, 10).Value = 7136.2856
$ws.Cells.Item(41, 11).Value = 11979.4
$ws.Cells.Item(41, 12).Value = 7136.2856
$ws.Cells.Item(41, 13).Value = -11589.4
$ws.Cells.Item(41, 14).Value = -7916.2856

$ws.Cells.Item(132, 8).Value = 2717.3845
$ws.Cells.Item(132, 9).Value = 2592.7646
$ws.Cells.Item(132, 10).Value = 2952.7778
$ws.Cells.Item(132, 11).Value = 7778.293799999999
$ws.Cells.Item(132, 12).Value = 8858.3334
$ws.Cells.Item(132, 13).Value = -5248.293799999999
$ws.Cells.Item(132, 14).Value = -13918.3334
